$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows for the accounts that were removed:
#   004996634 / HIROKO   / 22533.2  (row 2)
#   004376853 / ALBERTO  / 21092.39 (row 3)
#   004322549 / SIMONE   / 1153.67  (row 5, after prior deletions shifts to row 4)
# Row 005064129 / THIAGO / 20354.42 must remain in place.

$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(3).Delete()
